$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "name" column (C) to hold the
# new "prolificid" field. This shifts name/gender/race/re_rank one
# column to the right and keeps their existing formats (and, for most
# rows, their existing values).
$ws.Columns("C:C").Insert()

# ----- Header row -----
$ws.Range("C1").Value = "prolificid"

# Add the new "re_rank" header in column H, matching the look of the
# other header cells (bold font + border, same as G1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("H1").Value = "re_rank"

# ----- Data rows -----
# Columns: A=index, B=offer#, C=prolificid, D=name, E=gender, F=realeffort, G=race, H=re_rank
# Row 2 (index 0) - Giana
$ws.Range("B2").Value = 41
$ws.Range("C2").Value = "60bfcf5805c5ae12a546f9f3"
$ws.Range("D2").Value = "Giana"
$ws.Range("E2").Value = "female"
$ws.Range("F2").Value = 7.300264937320475
$ws.Range("G2").Value = "White"
$ws.Range("H2").Value = 1

# Row 3 (index 1) - Jewel
$ws.Range("B3").Value = 19
$ws.Range("C3").Value = "60b45e9961dd412bfb6780f8"
$ws.Range("D3").Value = "Jewel"
$ws.Range("E3").Value = "female"
$ws.Range("F3").Value = 6.476670993744667
$ws.Range("G3").Value = "Black or African American"
$ws.Range("H3").Value = 2

# Row 4 (index 2) - Colleen
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "5c5882fc5bfe7600011197cb"
$ws.Range("D4").Value = "Colleen"
$ws.Range("E4").Value = "female"
$ws.Range("F4").Value = 6.369967191149581
$ws.Range("G4").Value = "White"
$ws.Range("H4").Value = 3

# Row 5 (index 3) - Annes
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "60bd88b8fc436774352f53b9"
$ws.Range("D5").Value = "Annes"
$ws.Range("E5").Value = "female"
$ws.Range("F5").Value = 5.068061057111064
$ws.Range("G5").Value = "Asian"
$ws.Range("H5").Value = 4

# Row 6 (index 4) - Tina
$ws.Range("B6").Value = 34
$ws.Range("C6").Value = "5e96194b0a9fe909389e9f7b"
$ws.Range("D6").Value = "Tina"
$ws.Range("E6").Value = "female"
$ws.Range("F6").Value = 5.020468814187423
$ws.Range("G6").Value = "White"
$ws.Range("H6").Value = 5

# Row 7 (index 5) - Nansi
$ws.Range("B7").Value = 44
$ws.Range("C7").Value = "60c0e5899d387663c07eb3a4"
$ws.Range("D7").Value = "Nansi"
$ws.Range("E7").Value = "female"
$ws.Range("F7").Value = 4.300836668514529
$ws.Range("G7").Value = "Asian"
$ws.Range("H7").Value = 6

# Row 8 (index 6) - Lori (swapped with Khushi vs. the original order)
$ws.Range("B8").Value = 35
$ws.Range("C8").Value = "6077db0613ce87b4a62a78f9"
$ws.Range("D8").Value = "Lori"
$ws.Range("E8").Value = "female"
$ws.Range("F8").Value = 1.243334033638253
$ws.Range("G8").Value = "White"
$ws.Range("H8").Value = 7

# Row 9 (index 7) - Khushi (swapped with Lori vs. the original order)
$ws.Range("B9").Value = 22
$ws.Range("C9").Value = "608b14a312c099ac00b721b6"
$ws.Range("D9").Value = "Khushi"
$ws.Range("E9").Value = "female"
$ws.Range("F9").Value = 1.083571972669488
$ws.Range("G9").Value = "Asian"
$ws.Range("H9").Value = 8

# Row 10 (index 8) - Shaniek
$ws.Range("B10").Value = 33
$ws.Range("C10").Value = "60cb36ee9f58331a33cf5506"
$ws.Range("D10").Value = "Shaniek"
$ws.Range("E10").Value = "female"
$ws.Range("F10").Value = 0.4793892001648432
$ws.Range("G10").Value = "Black or African American"
$ws.Range("H10").Value = 9

# Row 11 (index 9) - Bri
$ws.Range("B11").Value = 21
$ws.Range("C11").Value = "5c0e89c6c323400001e6c4a5"
$ws.Range("D11").Value = "Bri"
$ws.Range("E11").Value = "female"
$ws.Range("F11").Value = 0.4192340444739328
$ws.Range("G11").Value = "Black or African American"
$ws.Range("H11").Value = 10

# Row 12 (index 10) - Kellie
$ws.Range("B12").Value = 32
$ws.Range("C12").Value = "6036f9b3b1842f8b659b18c7"
$ws.Range("D12").Value = "Kellie"
$ws.Range("E12").Value = "female"
$ws.Range("F12").Value = 0.2565399071127668
$ws.Range("G12").Value = "White"
$ws.Range("H12").Value = 11

# Row 13 (index 11) - Shadaisia
$ws.Range("B13").Value = 30
$ws.Range("C13").Value = "60d5775a99b502eec8cf56b4"
$ws.Range("D13").Value = "Shadaisia"
$ws.Range("E13").Value = "female"
$ws.Range("F13").Value = 0.2319910726680612
$ws.Range("G13").Value = "Black or African American"
$ws.Range("H13").Value = 12
